$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 25 ---
$ws.Range("B25").Value = 99874

# --- Row 26 ---
$ws.Range("B26").Value = 94292

# --- Row 27 (becomes the "Vårärt" / Lathyrus vernus record) ---
$ws.Range("A27").Value = 112395263
$ws.Range("B27").Value = 99874
$ws.Range("E27").Value = 221235
$ws.Range("F27").Value = "Vårärt"
$ws.Range("G27").Value = "Lathyrus vernus"
$ws.Range("H27").Value = "(L.) Bernh."
$ws.Range("J27").ClearContents()
$ws.Range("K27").ClearContents()
$ws.Range("L27").ClearContents()
$ws.Range("M27").ClearContents()
$ws.Range("N27").ClearContents()
$ws.Range("Q27").Value = 331849
$ws.Range("R27").Value = 6626616
$ws.Range("AF27").ClearContents()

# --- Row 28 ---
$ws.Range("A28").Value = 112395267
$ws.Range("B28").Value = 77650
$ws.Range("Q28").Value = 331734
$ws.Range("R28").Value = 6626659

# --- Row 29 (becomes the "Kattfotslav" / Felipes leucopellaeus record) ---
$ws.Range("A29").Value = 112395269
$ws.Range("B29").Value = 73772
$ws.Range("E29").Value = 6426
$ws.Range("F29").Value = "Kattfotslav"
$ws.Range("G29").Value = "Felipes leucopellaeus"
$ws.Range("H29").Value = "(Ach.) Frisch & G.Thor"
$ws.Range("Q29").Value = 331242
$ws.Range("R29").Value = 6626564

# --- Row 30 (becomes the "Vågbandad barkbock" / Semanotus undatus record) ---
$ws.Range("A30").Value = 112395265
$ws.Range("B30").Value = 5135
$ws.Range("E30").Value = 105930
$ws.Range("F30").Value = "Vågbandad barkbock"
$ws.Range("G30").Value = "Semanotus undatus"
$ws.Range("H30").Value = "(Linnaeus, 1758)"
$ws.Range("J30").Value = ""
$ws.Range("K30").Value = ""
$ws.Range("L30").Value = ""
$ws.Range("M30").Value = "äldre gnagspår"
$ws.Range("N30").Value = ""
$ws.Range("Q30").Value = 331818
$ws.Range("R30").Value = 6626574
$ws.Range("AF30").Value = ""

# --- Row 31 ---
$ws.Range("A31").Value = 112395266
$ws.Range("B31").Value = 77650
$ws.Range("Q31").Value = 331766
$ws.Range("R31").Value = 6626669

# --- Row 32 ---
$ws.Range("B32").Value = 90814
